$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Cntn2", "Nrcam", "ECs", 1, 0.5, 0.0700585, 0.140117, 0.4535220567529688, 0.4377808049690373, 2, 1, 1.172497, 2.344994, 0.04941711537896878, 0.04793039107701846, 0.08214338107449999, 0.328573524298, 0.02241175180546869, 0.0209830051881779),
  @("ECs", "Cntn2", "Nrcam", "FAPs", 1, 0.5, 0.0700585, 0.140117, 0.4535220567529688, 0.4377808049690373, 3, 1, 1.423194333333333, 4.269583, 0.05998323115285403, 0.08726793455581965, 0.09970686020183332, 0.5982411612109999, 0.02720371836313111, 0.03820422663783199),
  @("ECs", "Cntn2", "Nrcam", "Inflammatory-Mac", 1, 0.5, 0.0700585, 0.140117, 0.4535220567529688, 0.4377808049690373, 1, 0.3333333333333333, 0.02188533333333334, 0.06565600000000001, 0.0009223989847654408, 0.001341972626178457, 0.001533253625333334, 0.009199521752000001, 0.000418328284717673, 0.0005874898565348179),
  @("ECs", "Cntn2", "Nrcam", "MuSCs", 1, 0.5, 0.0700585, 0.140117, 0.4535220567529688, 0.4377808049690373, 2, 1, 21.082121, 42.164242, 0.8885460737983815, 0.8618139784264042, 1.4769817740785, 5.907927096314, 0.4029752429088171, 0.3772856172090797),
  @("ECs", "Cntn2", "Nrcam", "Neutrophils", 1, 0.5, 0.0700585, 0.140117, 0.4535220567529688, 0.4377808049690373, 1, 0.3333333333333333, 0.004358333333333333, 0.013075, 0.0001836902450013424, 0.0002672458280626801, 0.0003053382958333333, 0.001832029775, 0.00008330757771846554, 0.000116995093733897),
  @("ECs", "Cntn2", "Nrcam", "Resolving-Mac", 1, 0.5, 0.0700585, 0.140117, 0.4535220567529688, 0.4377808049690373, 2, 0.6666666666666666, 0.02248066666666667, 0.067442, 0.0009474904400291039, 0.001378477486516503, 0.001574961785666666, 0.009449770713999999, 0.0004297078131157747, 0.0006034709836788897),
  @("Inflammatory-Mac", "Cntn2", "Nrcam", "ECs", 1, 0.3333333333333333, 0.011109, 0.033327, 0.07191385097409637, 0.1041267004517875, 2, 1, 1.172497, 2.344994, 0.04941711537896878, 0.04793039107701846, 0.013025269173, 0.078151615038, 0.003553775070932887, 0.004990833474213729),
  @("Inflammatory-Mac", "Cntn2", "Nrcam", "FAPs", 1, 0.3333333333333333, 0.011109, 0.033327, 0.07191385097409637, 0.1041267004517875, 3, 1, 1.423194333333333, 4.269583, 0.05998323115285403, 0.08726793455581965, 0.015810265849, 0.142292392641, 0.004313625146071119, 0.009086922080540027),
  @("Inflammatory-Mac", "Cntn2", "Nrcam", "Inflammatory-Mac", 1, 0.3333333333333333, 0.011109, 0.033327, 0.07191385097409637, 0.1041267004517875, 1, 0.3333333333333333, 0.02188533333333334, 0.06565600000000001, 0.0009223989847654408, 0.001341972626178457, 0.0002431241680000001, 0.002188117512, 0.0000663332631290797, 0.0001397351816605828),
  @("Inflammatory-Mac", "Cntn2", "Nrcam", "MuSCs", 1, 0.3333333333333333, 0.011109, 0.033327, 0.07191385097409637, 0.1041267004517875, 2, 1, 21.082121, 42.164242, 0.8885460737983815, 0.8618139784264042, 0.234201282189, 1.405207693134, 0.06389876993475524, 0.08973784597676943),
  @("Inflammatory-Mac", "Cntn2", "Nrcam", "Neutrophils", 1, 0.3333333333333333, 0.011109, 0.033327, 0.07191385097409637, 0.1041267004517875, 1, 0.3333333333333333, 0.004358333333333333, 0.013075, 0.0001836902450013424, 0.0002672458280626801, 0.000048416725, 0.000435750525, 0.00001320987290442179, 0.00002782742628567259),
  @("Inflammatory-Mac", "Cntn2", "Nrcam", "Resolving-Mac", 1, 0.3333333333333333, 0.011109, 0.033327, 0.07191385097409637, 0.1041267004517875, 2, 0.6666666666666666, 0.02248066666666667, 0.067442, 0.0009474904400291039, 0.001378477486516503, 0.000249737726, 0.002247639534, 0.00006813768630363397, 0.0001435363123180368),
  @("MuSCs", "Cntn2", "Nrcam", "ECs", 1, 0.5, 0.073309, 0.146618, 0.4745640922729347, 0.4580924945791753, 2, 1, 1.172497, 2.344994, 0.04941711537896878, 0.04793039107701846, 0.085954582573, 0.343818330292, 0.0234515885025672, 0.02195655241462683),
  @("MuSCs", "Cntn2", "Nrcam", "FAPs", 1, 0.5, 0.073309, 0.146618, 0.4745640922729347, 0.4580924945791753, 3, 1, 1.423194333333333, 4.269583, 0.05998323115285403, 0.08726793455581965, 0.1043329533823333, 0.625997720294, 0.02846588764365179, 0.03997678583744764),
  @("MuSCs", "Cntn2", "Nrcam", "Inflammatory-Mac", 1, 0.5, 0.073309, 0.146618, 0.4745640922729347, 0.4580924945791753, 1, 0.3333333333333333, 0.02188533333333334, 0.06565600000000001, 0.0009223989847654408, 0.001341972626178457, 0.001604391901333334, 0.009626351408, 0.0004377374369186879, 0.0006147475879830566),
  @("MuSCs", "Cntn2", "Nrcam", "MuSCs", 1, 0.5, 0.073309, 0.146618, 0.4745640922729347, 0.4580924945791753, 2, 1, 21.082121, 42.164242, 0.8885460737983815, 0.8618139784264042, 1.545509208389, 6.182036833556, 0.4216720609548089, 0.3947905152405551),
  @("MuSCs", "Cntn2", "Nrcam", "Neutrophils", 1, 0.5, 0.073309, 0.146618, 0.4745640922729347, 0.4580924945791753, 1, 0.3333333333333333, 0.004358333333333333, 0.013075, 0.0001836902450013424, 0.0002672458280626801, 0.0003195050583333333, 0.00191703035, 0.00008717279437845502, 0.0001224233080431105),
  @("MuSCs", "Cntn2", "Nrcam", "Resolving-Mac", 1, 0.5, 0.073309, 0.146618, 0.4745640922729347, 0.4580924945791753, 2, 0.6666666666666666, 0.02248066666666667, 0.067442, 0.0009474904400291039, 0.001378477486516503, 0.001648035192666667, 0.009888211156, 0.0004496449406096951, 0.0006314701905195762)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
